# Functional Requirements.xlsx - add 2 requirements / user stories #58-67
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing rows 63-66 with new requirement text & dates ---
# (processing order matches the order the new strings were first introduced
#  so the rebuilt shared-string table lines up with the authored edit)

$ws.Range("B63").Value = "Pupil's details can be viewed by any teacher teaching him"
$ws.Range("E63").Value = 42703

# --- New row 68 : "A list of all pupils in the system can be viewed" ---
$ws.Range("B66").Copy()
$ws.Range("B68").PasteSpecial(-4122)
$ws.Range("C66:E66").Copy()
$ws.Range("C68:E68").PasteSpecial(-4122)
$ws.Range("G66").Copy()
$ws.Range("G68").PasteSpecial(-4122)

$ws.Range("B68").Value = "A list of all pupils in the system can be viewed"
$ws.Range("C68").Value = 5
$ws.Range("D68").Value = 2
$ws.Range("E68").Value = 42690
$ws.Range("G68").Value = "New"

# --- New row 67 : "A list of al teachers in the system can be viewed" ---
$ws.Range("B66").Copy()
$ws.Range("B67").PasteSpecial(-4122)
$ws.Range("C66:E66").Copy()
$ws.Range("C67:E67").PasteSpecial(-4122)
$ws.Range("G66").Copy()
$ws.Range("G67").PasteSpecial(-4122)

$ws.Range("B67").Value = "A list of al teachers in the system can be viewed"
$ws.Range("C67").Value = 5
$ws.Range("D67").Value = 2
$ws.Range("E67").Value = 42690
$ws.Range("G67").Value = "New"

$ws.Range("B64").Value = "Pupil's details can be viewed by himself"
$ws.Range("E64").Value = 42703

$ws.Range("B65").Value = "Details about a teacher can be viewed by all teachers"
$ws.Range("E65").Value = 42703

$ws.Range("B66").Value = "Details about a teacher can be viewed by their pupils"
$ws.Range("E66").Value = 42703

# --- New empty row 72 (extends used range to G72), matching row71's style ---
$ws.Range("A71").Copy()
$ws.Range("A72").PasteSpecial(-4122)

# --- Update the view: selection moves to B66 ---
$ws.Range("B66").Select()
